$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so values
# such as "44.80" or "0.00001047" are not auto-converted to numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.190.47'
$ws.Range('E2').Value = '  -3.73%  '

$ws.Range('D3').Value = '1.747.31'
$ws.Range('E3').Value = '  -2.83%  '

$ws.Range('D4').Value = '1.018'
$ws.Range('E4').Value = '  +1.55%  '

$ws.Range('D5').Value = '322.91'
$ws.Range('E5').Value = '  -4.70%  '

$ws.Range('D6').Value = '1.014'
$ws.Range('E6').Value = '  +1.48%  '

$ws.Range('D7').Value = '0.4196'
$ws.Range('E7').Value = '  -12.30%  '

$ws.Range('D8').Value = '0.3563'
$ws.Range('E8').Value = '  -2.19%  '

$ws.Range('D9').Value = '44.80'
$ws.Range('E9').Value = '  -1.41%  '

$ws.Range('D10').Value = '1.107'
$ws.Range('E10').Value = '  -3.22%  '

$ws.Range('D11').Value = '0.07321'
$ws.Range('E11').Value = '  -4.75%  '

$ws.Range('D12').Value = '1.017'
$ws.Range('E12').Value = '  +1.61%  '

$ws.Range('D13').Value = '21.43'
$ws.Range('E13').Value = '  -5.13%  '

$ws.Range('D14').Value = '6.054'
$ws.Range('E14').Value = '  -3.87%  '

$ws.Range('D15').Value = '7.160'
$ws.Range('E15').Value = '  -1.93%  '

$ws.Range('D16').Value = '1.751.02'
$ws.Range('E16').Value = '  -2.45%  '

$ws.Range('D17').Value = '0.00001047'
$ws.Range('E17').Value = '  -4.19%  '

$ws.Range('D18').Value = '83.35'
$ws.Range('E18').Value = '  +2.01%  '

$ws.Range('D19').Value = '1.013'
$ws.Range('E19').Value = '  +1.41%  '

$ws.Range('D20').Value = '0.05857'
$ws.Range('E20').Value = '  -12.78%  '

$ws.Range('D21').Value = '16.65'
$ws.Range('E21').Value = '  -3.77%  '

$ws.Range('D22').Value = '6.054'
$ws.Range('E22').Value = '  -5.43%  '

$ws.Range('D23').Value = '27.281.34'
$ws.Range('E23').Value = '  -3.39%  '

$ws.Range('D24').Value = '11.16'
$ws.Range('E24').Value = '  -6.98%  '

$ws.Range('D25').Value = '2.421'
$ws.Range('E25').Value = '  +0.88%  '

$ws.Range('D26').Value = '19.80'
$ws.Range('E26').Value = '  -3.71%  '

$ws.Range('D27').Value = '148.56'
$ws.Range('E27').Value = '  -1.86%  '

$ws.Range('D28').Value = '2.298'
$ws.Range('E28').Value = '  -4.35%  '

$ws.Range('D29').Value = '1.954.45'
$ws.Range('E29').Value = '  -2.41%  '

$ws.Range('D30').Value = '1.233'
$ws.Range('E30').Value = '  -2.89%  '

$ws.Range('D31').Value = '126.22'
$ws.Range('E31').Value = '  -5.61%  '

$ws.Range('D32').Value = '3.670'
$ws.Range('E32').Value = '  -9.13%  '

$ws.Range('D33').Value = '0.08994'
$ws.Range('E33').Value = '  -6.95%  '

$ws.Range('D34').Value = '5.498'
$ws.Range('E34').Value = '  -7.13%  '

$ws.Range('D35').Value = '12.36'
$ws.Range('E35').Value = '  +1.71%  '

$ws.Range('D36').Value = '0.2128'
$ws.Range('E36').Value = '  -2.04%  '

$ws.Range('D37').Value = '0.06052'
$ws.Range('E37').Value = '  -3.79%  '

$ws.Range('D38').Value = '0.02230'
$ws.Range('E38').Value = '  -6.39%  '

$ws.Range('D39').Value = '0.6355'
$ws.Range('E39').Value = '  -4.54%  '

$ws.Range('D40').Value = '4.939'
$ws.Range('E40').Value = '  -5.21%  '

$ws.Range('B41').Value = 'Frax'
$ws.Range('C41').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D41').Value = '1.013'
$ws.Range('E41').Value = '  +1.45%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '1.170'
$ws.Range('E42').Value = '  -3.76%  '

$ws.Range('B43').Value = 'WEMIXTOKEN'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '1.420'
$ws.Range('E43').Value = '  -4.16%  '

$ws.Range('D44').Value = '7.897'
$ws.Range('E44').Value = '  -2.32%  '

$ws.Range('D45').Value = '13.47'
$ws.Range('E45').Value = '  -4.37%  '

$ws.Range('D46').Value = '3.754'
$ws.Range('E46').Value = '  -2.88%  '

$ws.Range('D47').Value = '0.5805'
$ws.Range('E47').Value = '  -5.22%  '

$ws.Range('D48').Value = '122.78'
$ws.Range('E48').Value = '  -4.31%  '

$ws.Range('D49').Value = '1.926'
$ws.Range('E49').Value = '  -5.47%  '

$ws.Range('D50').Value = '0.06841'
$ws.Range('E50').Value = '  -3.62%  '

$ws.Range('D51').Value = '1.094'
$ws.Range('E51').Value = '  -6.65%  '
